$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (A1:D1) from numeric placeholders to descriptive text labels
$ws.Range("A1").Value = "Train RMSE"
$ws.Range("B1").Value = "Test RMSE"
$ws.Range("C1").Value = "Train MAE"
$ws.Range("D1").Value = "Test MAE"

# Update row 2 values (tiny floating point re-computation differences)
$ws.Range("A2").Value = 804.1864288928818
$ws.Range("B2").Value = 543.7011417795202
$ws.Range("C2").Value = 530.3719787638131
$ws.Range("D2").Value = 440.053699993206

# Remove now-duplicate rows 3 through 11 so only the header + one data row remain
$ws.Range("A3:D11").EntireRow.Delete()
